$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.415.56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.643.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "298.95"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3785"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.01"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3515"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08069"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.208"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.93%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.01"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.357"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.288"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001198"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.636.00"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.70"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06967"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.731"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.32"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.35"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.415.15"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.481"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.884"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.28%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.95"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.213"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.40"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.821.26"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.834"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.129"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.48"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9812"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -9.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02701"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08734"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2429"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.895"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06786"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.27%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6839"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.287"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.57"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.75%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6315"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.246"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.900"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07717"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "126.79"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.139"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.65%  "
